$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: introduce a typo into the login URL ("https" -> "h1ttps") ---
$ws.Cells.Item(2, 5).Value = "h1ttps://login.api.guxiansheng.cn/index.php?c=user&a=login"

# --- Row 3: keep the (correct) URL text, but give the cell its own distinct
#            formatting (a border), same visual family as the other data
#            cells in the row ---
$ws.Range("E3").Borders.LineStyle = 1

# --- Row 4: fill in a brand-new row that duplicates row 2's content, except
#            it keeps the *working* (un-typo'd) URL ---
$ws.Cells.Item(4, 1).Value = "登录"
$ws.Cells.Item(4, 2).Value = "00000001"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = "POST"
$ws.Cells.Item(4, 5).Value = "https://login.api.guxiansheng.cn/index.php?c=user&a=login"
$ws.Cells.Item(4, 5).Borders.LineStyle = 0
$ws.Cells.Item(4, 6).Value = "{""username"":""|mobile|"",""password"":""|password|""}"
$ws.Cells.Item(4, 6).WrapText = $true
$ws.Cells.Item(4, 7).Value = "{`n        ""Content-Type"":""application/x-www-form-urlencoded; charset=UTF-8""`n    }"
$ws.Cells.Item(4, 7).WrapText = $true
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 1
$ws.Rows.Item(4).RowHeight = 165

# --- Selection moves to E2 ---
$ws.Range("E2").Select()
